$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update loading_percent values for the 380 kV case (rows 2-25)
$data = @{
    2 = @{ "B"="7.778358071332719"; "D"="3.274679922277584"; "E"="24.61621265281816"; "F"="16.98616095398754"; "G"="3.567782193506756"; "M"="41.81992765739399"; "O"="14.83333191811389" }
    3 = @{ "B"="7.70548892209343"; "D"="3.187549142892391"; "E"="23.25989930604221"; "F"="17.01074258240347"; "G"="3.571214819775512"; "M"="39.30861327673671"; "O"="14.96678832354679" }
    4 = @{ "B"="7.662239339491035"; "D"="3.132281646194434"; "E"="22.39362088992859"; "F"="17.04123063946369"; "G"="3.573417636947942"; "M"="37.67728630445013"; "O"="15.06074526215347" }
    5 = @{ "B"="7.645008957279188"; "D"="3.109332478577533"; "E"="22.03261736265197"; "F"="17.05744966350543"; "G"="3.574339359313176"; "M"="36.99015830309371"; "O"="15.10198085689085" }
    6 = @{ "B"="7.642172211232761"; "D"="3.105496528683843"; "E"="21.97220495427291"; "F"="17.06036966057875"; "G"="3.574493867185133"; "M"="36.8747159739326"; "O"="15.10900386009095" }
    7 = @{ "B"="7.662005344066526"; "D"="3.131973851094981"; "E"="22.38878397556407"; "F"="17.04143411523772"; "G"="3.573429970048394"; "M"="37.66810974789654"; "O"="15.06128954313208" }
    8 = @{ "B"="7.752934595544235"; "D"="3.245014406932632"; "E"="24.15574339177918"; "F"="16.99139784393044"; "G"="3.568946086559295"; "M"="40.97263504328001"; "O"="14.87680981221339" }
    9 = @{ "B"="7.942226813206298"; "D"="3.451961005869444"; "E"="27.47505586749453"; "F"="17.0184414535519"; "G"="3.560902491409383"; "M"="46.74047786953941"; "O"="14.61365799427336" }
    10 = @{ "B"="8.086806787318137"; "D"="3.594214406359426"; "E"="29.92356489215911"; "F"="17.11810104542919"; "G"="3.555441239110157"; "M"="50.54204109767078"; "O"="14.48488722411444" }
    11 = @{ "B"="8.153517255638389"; "D"="3.656661203950581"; "E"="30.97669438597977"; "F"="17.18131239775872"; "G"="3.55305227136055"; "M"="52.17692866468987"; "O"="14.44125735241736" }
    12 = @{ "B"="8.178892179464901"; "D"="3.679972395777082"; "E"="31.36680688401657"; "F"="17.20784972027014"; "G"="3.552161198309221"; "M"="52.78246701895058"; "O"="14.4269595450458" }
    13 = @{ "B"="8.173422555426043"; "D"="3.674967010700367"; "E"="31.28317464563257"; "F"="17.20201840521342"; "G"="3.552352505227574"; "M"="52.65265593717303"; "O"="14.42993876942307" }
    14 = @{ "B"="8.155602725072333"; "D"="3.658585819435196"; "E"="31.00896294941251"; "F"="17.18344338526124"; "G"="3.552978690943989"; "M"="52.22701819505923"; "O"="14.44003607658352" }
    15 = @{ "B"="8.144701641432119"; "D"="3.648507817512785"; "E"="30.8398708081209"; "F"="17.17240503208285"; "G"="3.553364011754981"; "M"="51.96453813590502"; "O"="14.44651274687891" }
    16 = @{ "B"="8.082464132006413"; "D"="3.590086868979272"; "E"="29.85352247557864"; "F"="17.11433254623553"; "G"="3.555599271327337"; "M"="50.43329835577831"; "O"="14.48804574640871" }
    17 = @{ "B"="8.044508529547851"; "D"="3.553659606428386"; "E"="29.23291045253422"; "F"="17.0833100325424"; "G"="3.556994862786001"; "M"="49.4697514424146"; "O"="14.51740716551807" }
    18 = @{ "B"="8.022767299465093"; "D"="3.532494949848328"; "E"="28.87023336590379"; "F"="17.06714810704606"; "G"="3.557806556494204"; "M"="48.90665086532238"; "O"="14.53569779357206" }
    19 = @{ "B"="8.015422195133965"; "D"="3.525292773571887"; "E"="28.74645389438333"; "F"="17.06196365307111"; "G"="3.558082929638948"; "M"="48.71446662820038"; "O"="14.54212930236887" }
    20 = @{ "B"="8.048539831653814"; "D"="3.557559439213458"; "E"="29.29956709208617"; "F"="17.08643816242983"; "G"="3.556845370629845"; "M"="49.57324300130869"; "O"="14.51413587097255" }
    21 = @{ "B"="8.160833946848149"; "D"="3.663406568903989"; "E"="31.08974088106011"; "F"="17.18882855274881"; "G"="3.552794397659415"; "M"="52.35240600377331"; "O"="14.43700931728903" }
    22 = @{ "B"="8.234873823886492"; "D"="3.730621100634569"; "E"="32.20912678483639"; "F"="17.27091117489746"; "G"="3.550225924329538"; "M"="54.08974446931409"; "O"="14.3995977356637" }
    23 = @{ "B"="8.19530521654373"; "D"="3.694930113454535"; "E"="31.61630241533156"; "F"="17.22570734515022"; "G"="3.551589578462113"; "M"="53.16971109658362"; "O"="14.41835206512246" }
    24 = @{ "B"="8.046717029264325"; "D"="3.555797015921911"; "E"="29.26944994036599"; "F"="17.0850187252606"; "G"="3.55691292690109"; "M"="49.52648295276219"; "O"="14.51561043595065" }
    25 = @{ "B"="7.889963869738946"; "D"="3.397642713030805"; "E"="26.51927964489767"; "F"="16.99733063575851"; "G"="3.562999140272399"; "M"="45.25667924382275"; "O"="14.67382806441247" }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = [double]$rowVals[$col]
    }
}
